$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.025.54"
$ws.Range("E2").Value = "  +0.95%  "
$ws.Range("D3").Value = "'3.872.56"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'529.88"
$ws.Range("E5").Value = "  +8.72%  "
$ws.Range("D6").Value = "'143.19"
$ws.Range("E6").Value = "  -2.18%  "
$ws.Range("D7").Value = "'0.609"
$ws.Range("E7").Value = "  -2.28%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "'0.714"
$ws.Range("E9").Value = "  -3.97%  "
$ws.Range("D10").Value = "'0.170"
$ws.Range("E10").Value = "  -6.00%  "
$ws.Range("D11").Value = "'0.0000331"
$ws.Range("E11").Value = "  -6.25%  "
$ws.Range("D12").Value = "'41.87"
$ws.Range("E12").Value = "  -3.01%  "
$ws.Range("D13").Value = "'4.500.84"
$ws.Range("E13").Value = "  -0.44%  "
$ws.Range("D14").Value = "'10.17"
$ws.Range("E14").Value = "  -2.94%  "
$ws.Range("D15").Value = "'3.889.91"
$ws.Range("E15").Value = "  -0.26%  "
$ws.Range("D16").Value = "'13.91"
$ws.Range("E16").Value = "  -2.04%  "
$ws.Range("D17").Value = "'1.21"
$ws.Range("E17").Value = "  +6.00%  "
$ws.Range("E18").Value = "  -1.55%  "
$ws.Range("D19").Value = "'20.23"
$ws.Range("E19").Value = "  +1.04%  "
$ws.Range("D20").Value = "'69.012.85"
$ws.Range("E20").Value = "  +0.82%  "
$ws.Range("D21").Value = "'422.32"
$ws.Range("E21").Value = "  -2.14%  "
$ws.Range("D22").Value = "'3.38"
$ws.Range("E22").Value = "  -4.83%  "
$ws.Range("D23").Value = "'14.11"
$ws.Range("E23").Value = "  -4.31%  "
$ws.Range("D24").Value = "'87.25"
$ws.Range("E24").Value = "  -2.57%  "
$ws.Range("D25").Value = "'4.00"
$ws.Range("E25").Value = "  +7.23%  "
$ws.Range("D26").Value = "'11.34"
$ws.Range("E26").Value = "  -7.44%  "
$ws.Range("D27").Value = "'10.53"
$ws.Range("E27").Value = "  -3.76%  "
$ws.Range("D28").Value = "'36.15"
$ws.Range("E28").Value = "  -3.35%  "
$ws.Range("D29").Value = "'693.46"
$ws.Range("E29").Value = "  -2.75%  "
$ws.Range("D30").Value = "'13.13"
$ws.Range("E30").Value = "  -1.95%  "
$ws.Range("E31").Value = "  -4.19%  "
$ws.Range("D32").Value = "'2.84"
$ws.Range("E32").Value = "  -2.97%  "
$ws.Range("D33").Value = "'67.64"
$ws.Range("E33").Value = "  +9.39%  "
$ws.Range("D34").Value = "'0.433"
$ws.Range("E34").Value = "  +7.25%  "
$ws.Range("D35").Value = "'5.94"
$ws.Range("E35").Value = "  -2.21%  "
$ws.Range("D36").Value = "'0.0₃0851"
$ws.Range("E36").Value = "  -5.07%  "
$ws.Range("D37").Value = "'39.87"
$ws.Range("E37").Value = "  -2.24%  "
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  -0.37%  "
$ws.Range("D39").Value = "'0.147"
$ws.Range("E39").Value = "  -0.20%  "
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").Value = "'3.24"
$ws.Range("E41").Value = "  +4.90%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.0479"
$ws.Range("E42").Value = "  -3.50%  "
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "'3.17"
$ws.Range("E43").Value = "  +4.23%  "
$ws.Range("D44").Value = "'2.77"
$ws.Range("E44").Value = "  -6.13%  "
$ws.Range("E45").Value = "  +1.55%  "
$ws.Range("E46").Value = "  -1.72%  "
$ws.Range("D47").Value = "'2.97"
$ws.Range("E47").Value = "  +5.52%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "'0.0₆0346"
$ws.Range("E48").Value = "  -7.66%  "
$ws.Range("D49").Value = "'2.749.56"
$ws.Range("E49").Value = "  +13.93%  "
$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D50").Value = "'0.000269"
$ws.Range("E50").Value = "  +9.32%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "'144.61"
$ws.Range("E51").Value = "  +1.05%  "
